$d = $word.ActiveDocument

# The document currently holds exactly two paragraphs: an empty one,
# and one that carries the _GoBack bookmark. Insert a new "Heaps:"
# paragraph right in front of the bookmarked paragraph.
$bookmarkPara = $d.Paragraphs.Item(2)
$bookmarkPara.Range.InsertBefore("Heaps:`r")

# The bookmarked paragraph is now paragraph 3. Put the new bullet text
# at the very front of it, ahead of the bookmarks that already live
# there (this keeps the bookmarks inside the same paragraph).
$listPara = $d.Paragraphs.Item(3)
$insertion = $d.Range($listPara.Range.Start, $listPara.Range.Start)
$insertion.InsertBefore("clever way to insert a distance, x, y co-ordinate in heap is to insert all three of them in the heap")

# Turn that paragraph into a bulleted list item (List Paragraph style +
# a fresh bullet-list numbering definition, cloned from the document's
# existing Symbol-bullet list).
$listPara.Style = "ListParagraph"
$bulletGallery = $word.ListGalleries.Item(1)
$bulletTemplate = $bulletGallery.ListTemplates.Item(1)
$listPara.Range.ListFormat.ApplyListTemplateWithLevel($bulletTemplate, $false, 0, $false, 0)

# Add a new, empty bulleted list paragraph right after it. Paragraphs.Add
# copies the reference paragraph's formatting (style + list numbering),
# so the new paragraph automatically continues the same list (numId 4)
# instead of starting a brand new one.
$newListPara = $d.Paragraphs.Add($listPara.Range)
